$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Test Data")
$wsAlign = $wb.Worksheets.Item("Alignment Type")

# --- Alignment Type sheet: insert a new row 14 (Product Test 3 / Central / Intermediary 12) ---
$wsAlign.Rows.Item(14).Insert()
$wsAlign.Range("A14").Value = "Product Test 3"
$wsAlign.Range("B14").Value = "Central"
$wsAlign.Range("C14").Value = "Intermediary 12"

# Keep the autofilter range in sync with the grown data range
$wsAlign.AutoFilterMode = $false
$wsAlign.Range("A1:C21").AutoFilter() | Out-Null

# --- Test Data sheet: update correlation summary for Product Test 3 row (row 7) ---
$wsData.Range("H7").Value = "'2/4"
$wsData.Range("J7").Value = "Intermediary 12`nIntermediary 8"
$wsData.Range("J7").WrapText = $true

# --- defined name: _FilterDatabase now covers the extra row ---
$name = $wb.Names.Item("Alignment Type!_FilterDatabase")
$name.RefersTo = "='Alignment Type'!`$A`$1:`$C`$21"

# --- selections / active sheet ---
$wsAlign.Range("D14").Select() | Out-Null
$wsData.Activate() | Out-Null
$wsData.Range("J8").Select() | Out-Null
